# Add new column "Fertilizer Use Per Area" before the existing
# "Fertilizer Use Per Capita" column (column L), shifting all subsequent
# columns (L..R) one position to the right (M..S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at L; existing L:R data shifts to M:S.
$ws.Columns("L").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("L1").Value = "Fertilizer Use Per Area"

# New data values for the inserted column, by row.
$values = @{
    2  = -0.5600000000000001
    3  = -0.1
    4  = -0.45
    5  = 0.67
    6  = -0.58
    7  = -0.22
    8  = 0.15
    9  = 0.8100000000000001
    10 = 0.15
    12 = -0.07000000000000001
    15 = 0.66
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 12).Value = $values[$row]
}
